$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows (Índice, Distancia, max, min, Tempo)
$data = @(
    @(0, 9288.299999999999, 10036, 8426, 0.1565564076105753),
    @(1, 9425.533333333333, 10063, 8369, 0.1608503341674805),
    @(2, 9568, 10247, 8824, 0.1598768313725789),
    @(3, 9873.133333333333, 10544, 8754, 0.1571749130884806),
    @(4, 9114.1, 10237, 7908, 0.1541552464167277),
    @(5, 10040.93333333333, 10948, 9161, 0.1604896863301595),
    @(6, 9269.666666666666, 9893, 8459, 0.1620262384414673),
    @(7, 9117.966666666667, 10194, 8317, 0.1592001517613729),
    @(8, 8664.933333333332, 9498, 7300, 0.1553608735402425),
    @(9, 8893.633333333333, 9790, 8124, 0.1557687362035116)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
